# Applies the textual updates described by the diff:
# - updates the header date
# - updates the 25 multiplication problems in the table

$d = $word.ActiveDocument

$replacements = @(
    @{old = "2023-10-20 Friday"; new = "2023-10-21 Saturday"},
    @{old = "83×91="; new = "59×74="},
    @{old = "23×18="; new = "70×43="},
    @{old = "26×11="; new = "85×29="},
    @{old = "11×22="; new = "21×37="},
    @{old = "16×15="; new = "43×70="},
    @{old = "69×48="; new = "88×53="},
    @{old = "91×93="; new = "37×86="},
    @{old = "22×55="; new = "96×78="},
    @{old = "15×12="; new = "52×86="},
    @{old = "43×40="; new = "82×43="},
    @{old = "70×32="; new = "99×91="},
    @{old = "66×14="; new = "29×45="},
    @{old = "70×11="; new = "60×28="},
    @{old = "62×65="; new = "40×39="},
    @{old = "12×64="; new = "63×89="},
    @{old = "52×71="; new = "81×88="},
    @{old = "11×78="; new = "64×97="},
    @{old = "94×97="; new = "64×98="},
    @{old = "66×91="; new = "39×71="},
    @{old = "69×57="; new = "17×59="},
    @{old = "56×48="; new = "37×91="},
    @{old = "20×81="; new = "26×86="},
    @{old = "34×68="; new = "79×67="},
    @{old = "39×63="; new = "92×17="},
    @{old = "97×90="; new = "13×17="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
